$d = $word.ActiveDocument

# Anchor: a zero-length Range collapsed to the very START of the document's
# final (empty, paraId 19FDC771) paragraph. Inserting a full set of
# <w:p>...</w:p> OOXML fragments there via InsertXML splits cleanly into new
# sibling paragraphs placed *before* that final paragraph (collapsing at the
# END of the previous paragraph instead merges the fragment's runs into it,
# which is not what we want).
$lastPara = $d.Paragraphs.Last
$insPoint = $lastPara.Range
$insPoint.Collapse(1)

$fragment = '<w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t>30th</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t>Nov</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t xml:space="preserve"> – </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t>4</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:vertAlign w:val="superscript"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t>th</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t xml:space="preserve"> Dec</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:tab/>
        <w:t xml:space="preserve">What I did last week: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="708" w:firstLine="708"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t>More</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t xml:space="preserve">Development on the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t>Canbus</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t>. Reading more documentation.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="708"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t xml:space="preserve">What I am doing this week: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="708"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t xml:space="preserve">Reading more </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t xml:space="preserve">into My SQL and AWS, trying to work with </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t>MCUxpresso</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t xml:space="preserve"> as hardware has not arrived</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:tab/>
        <w:t xml:space="preserve">Are there any blockers: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t xml:space="preserve">Yes – hardware has not arrived </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-IE"/>
        </w:rPr>
      </w:pPr>
    </w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $fragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insPoint.InsertXML($xml) | Out-Null

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
